# Daily attendance processing - 2025-11-24 11:47:50
# Applies the updated attendance / coverage figures and "Recorded by" list
# reorderings for the Urogenital Y3 B25/26 session analysis workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a literal text value into a cell without Excel's
# "smart" auto-conversion turning percentage-looking strings (e.g.
# "31.6%") into real numeric percentages. We briefly mark the cell as
# Text, assign the literal string, then leave it (Text format renders
# the literal string exactly, same as the source data).
# ---------------------------------------------------------------------
function Set-TextValue($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# ---------------------------------------------------------------------
# "Recorded By" email list re-orderings (same people, new order) plus
# a few genuinely new distribution lists below.
# ---------------------------------------------------------------------

$emailsUro = "Veronia.rafat@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, heba@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G2").Value = $emailsUro
$ws.Range("G21").Value = $emailsUro
$ws.Range("G40").Value = $emailsUro

$emailsHisto = "Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg"
$ws.Range("G6").Value = $emailsHisto
$ws.Range("G44").Value = $emailsHisto

$emailsPhysioA = "yasmin.m.senosy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, marinasorial@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G19").Value = $emailsPhysioA
$ws.Range("G76").Value = $emailsPhysioA
$ws.Range("G95").Value = $emailsPhysioA

$emailsPhysioB = "yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G20").Value = $emailsPhysioB
$ws.Range("G38").Value = $emailsPhysioB
$ws.Range("G39").Value = $emailsPhysioB
$ws.Range("G57").Value = $emailsPhysioB
$ws.Range("G58").Value = $emailsPhysioB
$ws.Range("G77").Value = $emailsPhysioB
$ws.Range("G96").Value = $emailsPhysioB
$ws.Range("G115").Value = $emailsPhysioB

$emailsMicroA2 = "yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, dina.adel@med.asu.edu.eg"
$ws.Range("G28").Value = $emailsMicroA2
$ws.Range("G100").Value = $emailsMicroA2

$emailsCBL = "NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg"
$ws.Range("G43").Value = $emailsCBL

$emailsPharm = "Mohammedeltanany@med.asu.edu.eg, heba@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg"
$ws.Range("G59").Value = $emailsPharm
$ws.Range("G78").Value = $emailsPharm
$ws.Range("G97").Value = $emailsPharm

# ---------------------------------------------------------------------
# Class statistics block (K6:L10 metrics)
# ---------------------------------------------------------------------
$ws.Range("L6").Value = 36
$ws.Range("L7").Value = 2
$ws.Range("L8").Value = 76
Set-TextValue "L9" "31.6%"
Set-TextValue "L10" "42.7%"

# ---------------------------------------------------------------------
# Group coverage summary rows (row 16 = A2, row 18 = B1, row 19 = B2)
# ---------------------------------------------------------------------
$ws.Range("O16").Value = 7
$ws.Range("P16").Value = 0
Set-TextValue "R16" "36.8%"
Set-TextValue "S16" "44.9%"

$ws.Range("O18").Value = 6
$ws.Range("Q18").Value = 12
Set-TextValue "R18" "31.6%"
Set-TextValue "S18" "23.0%"

$ws.Range("O19").Value = 6
$ws.Range("P19").Value = 0
Set-TextValue "R19" "31.6%"
Set-TextValue "S19" "36.6%"

# ---------------------------------------------------------------------
# Sessions that moved from Not Recorded / Pending to Recorded.
# Re-use the existing "Recorded" look (green fill) by pasting the
# formatting from an already-Recorded row, then fill in the new data.
# ---------------------------------------------------------------------

# Row 24: Year 3 / A2 / BIOCHEMISTRY LAB/CBL -> now Recorded
$ws.Range("A2:I2").Copy() | Out-Null
$ws.Range("A24:I24").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("G24").Value = "Fatmaelhady@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg"
$ws.Range("H24").Value = "15/204"
$ws.Range("I24").Value = "Recorded"

# Row 66: Year 3 / B1 / MICROBIOLOGY -> now Recorded
$ws.Range("A2:I2").Copy() | Out-Null
$ws.Range("A66:I66").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("G66").Value = "amira.m.ibrahim@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
$ws.Range("H66").Value = "2/149"
$ws.Range("I66").Value = "Recorded"

# Row 85: Year 3 / B2 / MICROBIOLOGY -> now Recorded
$ws.Range("A2:I2").Copy() | Out-Null
$ws.Range("A85:I85").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("G85").Value = "amira.m.ibrahim@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
$ws.Range("H85").Value = "97/227"
$ws.Range("I85").Value = "Recorded"

Write-Output "done"
